$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 22
$ws.Range("F6").Value = 4150
$ws.Range("F7").Value = 1844
$ws.Range("A8").Value = 9
$ws.Range("A9").Value = 33
$ws.Range("F9").Value = 1031
$ws.Range("A10").Value = 25
$ws.Range("A11").Value = 0
$ws.Range("A15").Value = 42
$ws.Range("A16").Value = 24
$ws.Range("F16").Value = 2750
$ws.Range("A18").Value = 30
$ws.Range("A19").Value = 23
$ws.Range("A20").Value = 45
$ws.Range("A21").Value = 39
$ws.Range("C21").Value = "Fogo de Chão Brazilian Steakhouse"
$ws.Range("D21").Value = 3
$ws.Range("F21").Value = 4874
$ws.Range("A22").Value = 16
$ws.Range("C22").Value = "Greek Islands"
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 4.6
$ws.Range("F22").Value = 866
$ws.Range("A26").Value = 21
$ws.Range("A27").Value = 19
$ws.Range("F28").Value = 2791
$ws.Range("A30").Value = 8
$ws.Range("F30").Value = 2257
$ws.Range("A32").Value = 43
$ws.Range("A33").Value = 7
$ws.Range("A34").Value = 38
$ws.Range("F34").Value = 666
$ws.Range("A35").Value = 36
$ws.Range("A36").Value = 31
$ws.Range("C36").Value = "Ocean Prime"
$ws.Range("D36").Value = 4
$ws.Range("E36").Value = 4.6
$ws.Range("F36").Value = 958
$ws.Range("A37").Value = 59
$ws.Range("C37").Value = "Olive Garden Italian Restaurant"
$ws.Range("D37").Value = 2
$ws.Range("E37").Value = 4.2
$ws.Range("F37").Value = 2759
$ws.Range("A38").Value = 3
$ws.Range("C38").Value = "Rusty Bucket Restaurant and Tavern"
$ws.Range("F38").Value = 946
$ws.Range("A39").Value = 32
$ws.Range("C39").Value = "Ruth's Chris Steak House"
$ws.Range("D39").Value = 4
$ws.Range("E39").Value = 4.6
$ws.Range("F39").Value = 1709
$ws.Range("A40").Value = 2
$ws.Range("C40").Value = "Seasons 52"
$ws.Range("D40").Value = 2
$ws.Range("E40").Value = 4.5
$ws.Range("F40").Value = 1339
$ws.Range("A41").Value = 11
$ws.Range("C41").Value = "Sero's Family Restaurant"
$ws.Range("F41").Value = 1158
$ws.Range("A42").Value = 37
$ws.Range("C42").Value = "Slapfish"
$ws.Range("E42").Value = 4.6
$ws.Range("F42").Value = 317
$ws.Range("A43").Value = 40
$ws.Range("C43").Value = "Steer-In"
$ws.Range("D43").Value = 1
$ws.Range("E43").Value = 4.5
$ws.Range("F43").Value = 2703
$ws.Range("A44").Value = 35
$ws.Range("A46").Value = 1
$ws.Range("F46").Value = 821
$ws.Range("A47").Value = 29
$ws.Range("A48").Value = 20
$ws.Range("A49").Value = 58
$ws.Range("A50").Value = 26
$ws.Range("F52").Value = 3906
$ws.Range("A53").Value = 28
$ws.Range("A55").Value = 10
$ws.Range("A56").Value = 34
